$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B7: clarify clients task now includes "detalle"
$ws.Range("B7").Value = "Listado, detalle, creacion y eliminacion de clientes"

# Row 9: new task 6.3 - Productos listing/detail/create/delete
$ws.Range("B9").Value = "Listado, detalle, creacion y eliminacion de Productos "

# "6.3" must stay a text value (like "6.2" in A8), not become numeric 6.3,
# so format the cell as Text before typing it, then restore the original
# look (font/fill/border/alignment) by copying the formatting from A8.
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = "6.3"
$ws.Range("A8").Copy()
$ws.Range("A9").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("C9").Value = "Jon"
$ws.Range("D9").Value = [datetime]"2025-04-05"
$ws.Range("E9").Value = [datetime]"2025-04-05"
$ws.Range("F9").Value = "✅ Hecho"
$ws.Range("A9:G9").RowHeight = 30

# Row 10: new task 7 - Actualizacion del ReadMe
$ws.Range("A10").Value = 7
$ws.Range("B10").Value = "Actualizacion del ReadMe"
$ws.Range("C10").Value = "Jon"
$ws.Range("D10").Value = [datetime]"2025-04-05"
$ws.Range("E10").Value = [datetime]"2025-04-05"
$ws.Range("F10").Value = "✅ Hecho"

$excel.CutCopyMode = 0

# Update selected cell to reflect latest edit location
$ws.Range("H10").Select()
